$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = [double]"240.0054887063018"
$ws.Range("B2").Value = [double]"250.00804305218327"
$ws.Range("C2").Value = [double]"10.865066011000692"
$ws.Range("D2").Value = [double]"0.9999085846741057"
$ws.Range("E2").Value = [double]"-0.011537826972071044"
$ws.Range("F2").Value = [double]"0.007874446579951469"
$ws.Range("G2").Value = [double]"2.525960433510624e-5"

$ws.Range("A3").Value = [double]"240.00487883783265"
$ws.Range("B3").Value = [double]"249.9893275409291"
$ws.Range("C3").Value = [double]"11.650951820165472"
$ws.Range("D3").Value = [double]"0.9996213949776022"
$ws.Range("E3").Value = [double]"0.027912827853204314"
$ws.Range("F3").Value = [double]"-0.0009089997508487205"
$ws.Range("G3").Value = [double]"0.00026626139923763693"

$ws.Range("A4").Value = [double]"240.00198670229366"
$ws.Range("B4").Value = [double]"249.96938797230723"
$ws.Range("C4").Value = [double]"12.365117913307849"
$ws.Range("D4").Value = [double]"0.9995261248958033"
$ws.Range("E4").Value = [double]"0.030840655465116957"
$ws.Range("F4").Value = [double]"-0.0044640276236466185"
$ws.Range("G4").Value = [double]"0.00010096544427241551"

$ws.Range("A5").Value = [double]"239.99898295911967"
$ws.Range("B5").Value = [double]"249.94869246876027"
$ws.Range("C5").Value = [double]"13.013638531781247"
$ws.Range("D5").Value = [double]"0.9994525240047093"
$ws.Range("E5").Value = [double]"0.033134427452728356"
$ws.Range("F5").Value = [double]"-0.004797651442710175"
$ws.Range("G5").Value = [double]"0.00012511051918700198"

$ws.Range("A6").Value = [double]"239.99799229767285"
$ws.Range("B6").Value = [double]"249.92720740956142"
$ws.Range("C6").Value = [double]"13.602524261797527"
$ws.Range("D6").Value = [double]"0.9993816530357674"
$ws.Range("E6").Value = [double]"0.035532539240335816"
$ws.Range("F6").Value = [double]"-0.0016351250518979732"
$ws.Range("G6").Value = [double]"0.00014100492088642424"

$ws.Range("A7").Value = [double]"239.99695590253648"
$ws.Range("B7").Value = [double]"249.90479679737834"
$ws.Range("C7").Value = [double]"14.137422169546449"
$ws.Range("D7").Value = [double]"0.9992842944861353"
$ws.Range("E7").Value = [double]"0.038205248722148095"
$ws.Range("F7").Value = [double]"-0.001761736521247684"
$ws.Range("G7").Value = [double]"0.00019531652107795909"

$ws.Range("A8").Value = [double]"239.99574233179433"
$ws.Range("B8").Value = [double]"249.87806110631007"
$ws.Range("C8").Value = [double]"14.626896602905719"
$ws.Range("D8").Value = [double]"0.9989039595044177"
$ws.Range("E8").Value = [double]"0.04672192724976696"
$ws.Range("F8").Value = [double]"-0.0021207118123974435"
$ws.Range("G8").Value = [double]"0.0001809849440817684"

$ws.Range("A9").Value = [double]"239.9944827086574"
$ws.Range("B9").Value = [double]"249.8503944059455"
$ws.Range("C9").Value = [double]"15.071872967457303"
$ws.Range("D9").Value = [double]"0.9987625642925996"
$ws.Range("E9").Value = [double]"0.04963274824250498"
$ws.Range("F9").Value = [double]"-0.0022596500151978085"
$ws.Range("G9").Value = [double]"0.00025991317007576393"

$ws.Range("A10").Value = [double]"239.99322011181692"
$ws.Range("B10").Value = [double]"249.8223419570691"
$ws.Range("C10").Value = [double]"15.478936770321729"
$ws.Range("D10").Value = [double]"0.998671494826719"
$ws.Range("E10").Value = [double]"0.05149106127051035"
$ws.Range("F10").Value = [double]"-0.00231712963711739"
$ws.Range("G10").Value = [double]"1.7398599895875426e-6"

$ws.Range("A11").Value = [double]"239.99191418588873"
$ws.Range("B11").Value = [double]"249.7933272715605"
$ws.Range("C11").Value = [double]"15.849476529979126"
$ws.Range("D11").Value = [double]"0.9985130463367223"
$ws.Range("E11").Value = [double]"0.05447112987669662"
$ws.Range("F11").Value = [double]"-0.0024512567284888276"
$ws.Range("G11").Value = [double]"2.0500401159908975e-6"

$ws.Range("A12").Value = [double]"239.99056558436223"
$ws.Range("B12").Value = [double]"249.7633648735148"
$ws.Range("C12").Value = [double]"16.18688833353846"
$ws.Range("D12").Value = [double]"0.9983464527059074"
$ws.Range("E12").Value = [double]"0.057437270331758716"
$ws.Range("F12").Value = [double]"-0.0025847578850646744"
$ws.Range("G12").Value = [double]"2.3876884777263444e-6"

$ws.Range("A13").Value = [double]"239.98917546730908"
$ws.Range("B13").Value = [double]"249.73248051960547"
$ws.Range("C13").Value = [double]"16.494245123869984"
$ws.Range("D13").Value = [double]"0.9981736396045922"
$ws.Range("E13").Value = [double]"0.060359689775349876"
$ws.Range("F13").Value = [double]"-0.002716290786665449"
$ws.Range("G13").Value = [double]"2.747883340458671e-6"

$ws.Range("A14").Value = [double]"239.98774549585585"
$ws.Range("B14").Value = [double]"249.70071108349086"
$ws.Range("C14").Value = [double]"16.77432375057716"
$ws.Range("D14").Value = [double]"0.9979969005096729"
$ws.Range("E14").Value = [double]"0.06320826676096324"
$ws.Range("F14").Value = [double]"-0.002844498501955431"
$ws.Range("G14").Value = [double]"3.1286596356109124e-6"

$ws.Range("A15").Value = [double]"239.98627793569864"
$ws.Range("B15").Value = [double]"249.6681068416105"
$ws.Range("C15").Value = [double]"17.02962957007997"
$ws.Range("D15").Value = [double]"0.9978191865668616"
$ws.Range("E15").Value = [double]"0.06594799402570684"
$ws.Range("F15").Value = [double]"-0.002967803280911401"
$ws.Range("G15").Value = [double]"3.523403400136968e-6"

$ws.Range("A16").Value = [double]"239.98477569180505"
$ws.Range("B16").Value = [double]"249.63473224042232"
$ws.Range("C16").Value = [double]"17.262418495497723"
$ws.Range("D16").Value = [double]"0.9976439376291044"
$ws.Range("E16").Value = [double]"0.06854198020225062"
$ws.Range("F16").Value = [double]"-0.0030845421399806503"
$ws.Range("G16").Value = [double]"3.923951911600692e-6"

$ws.Range("A17").Value = [double]"239.98324235949997"
$ws.Range("B17").Value = [double]"249.60066702397984"
$ws.Range("C17").Value = [double]"17.474717241468266"
$ws.Range("D17").Value = [double]"0.9974751068870886"
$ws.Range("E17").Value = [double]"0.0709508086604734"
$ws.Range("F17").Value = [double]"-0.003192937792184588"
$ws.Range("G17").Value = [double]"4.318827265904371e-6"

$ws.Range("A18").Value = [double]"239.98168226497958"
$ws.Range("B18").Value = [double]"249.56600712953585"
$ws.Range("C18").Value = [double]"17.668342040434297"
$ws.Range("D18").Value = [double]"0.9973170906053108"
$ws.Range("E18").Value = [double]"0.07313304538652086"
$ws.Range("F18").Value = [double]"-0.0032911216716727887"
$ws.Range("G18").Value = [double]"4.6936839304746455e-6"

$ws.Range("A19").Value = [double]"239.98010050037965"
$ws.Range("B19").Value = [double]"249.53086547031714"
$ws.Range("C19").Value = [double]"17.844916122534052"
$ws.Range("D19").Value = [double]"0.9971746466450857"
$ws.Range("E19").Value = [double]"0.07504549146880408"
$ws.Range("F19").Value = [double]"-0.0033771458079456733"
$ws.Range("G19").Value = [double]"5.033561452105831e-6"

$ws.Range("A20").Value = [double]"239.9785029392296"
$ws.Range("B20").Value = [double]"249.49537226293694"
$ws.Range("C20").Value = [double]"18.005886290138253"
$ws.Range("D20").Value = [double]"0.9970527339248999"
$ws.Range("E20").Value = [double]"0.07664414031809665"
$ws.Range("F20").Value = [double]"-0.0034490248818666636"
$ws.Range("G20").Value = [double]"5.3150791817372455e-6"

$ws.Range("A21").Value = [double]"239.97689626180974"
$ws.Range("B21").Value = [double]"249.4596756124856"
$ws.Range("C21").Value = [double]"18.15253878913153"
$ws.Range("D21").Value = [double]"0.9969564317684267"
$ws.Range("E21").Value = [double]"0.07788359569851153"
$ws.Range("F21").Value = [double]"-0.0035047123703570823"
$ws.Range("G21").Value = [double]"5.519283048221191e-6"

$ws.Range("A22").Value = [double]"239.9757049871276"
$ws.Range("B22").Value = [double]"249.34399451386338"
$ws.Range("C22").Value = [double]"18.430681819097064"
$ws.Range("D22").Value = [double]"0.9736224203341091"
$ws.Range("E22").Value = [double]"0.22765302155570707"
$ws.Range("F22").Value = [double]"-0.002344597427024989"
$ws.Range("G22").Value = [double]"0.00037316161637441727"

$ws.Range("A23").Value = [double]"239.97441888688095"
$ws.Range("B23").Value = [double]"249.2269497937221"
$ws.Range("C23").Value = [double]"18.69275699886776"
$ws.Range("D23").Value = [double]"0.9724212378962032"
$ws.Range("E23").Value = [double]"0.2326823810134843"
$ws.Range("F23").Value = [double]"-0.002557151419175768"
$ws.Range("G23").Value = [double]"0.0006601345906706631"

$ws.Range("A24").Value = [double]"239.9730182700803"
$ws.Range("B24").Value = [double]"249.10847323378525"
$ws.Range("C24").Value = [double]"18.940802076311535"
$ws.Range("D24").Value = [double]"0.9712127728127671"
$ws.Range("E24").Value = [double]"0.23762547156879266"
$ws.Range("F24").Value = [double]"-0.002809785564950918"
$ws.Range("G24").Value = [double]"0.001024560084788972"

$ws.Range("A25").Value = [double]"239.97147185591095"
$ws.Range("B25").Value = [double]"248.98852037316112"
$ws.Range("C25").Value = [double]"19.176699935472733"
$ws.Range("D25").Value = [double]"0.9700092342130158"
$ws.Range("E25").Value = [double]"0.24243800846314537"
$ws.Range("F25").Value = [double]"-0.0031263089962544847"
$ws.Range("G25").Value = [double]"0.0015122141558332079"

$ws.Range("A26").Value = [double]"239.96977525246493"
$ws.Range("B26").Value = [double]"248.86710087296572"
$ws.Range("C26").Value = [double]"19.40213409938442"
$ws.Range("D26").Value = [double]"0.9688374419128001"
$ws.Range("E26").Value = [double]"0.2470249583317101"
$ws.Range("F26").Value = [double]"-0.0034527063769554697"
$ws.Range("G26").Value = [double]"0.0020187305537916"

$ws.Range("A27").Value = [double]"239.96764066853456"
$ws.Range("B27").Value = [double]"248.7473237799759"
$ws.Range("C27").Value = [double]"19.611184443404834"
$ws.Range("D27").Value = [double]"0.9689700102468797"
$ws.Range("E27").Value = [double]"0.24648074376270723"
$ws.Range("F27").Value = [double]"-0.004396083870233778"
$ws.Range("G27").Value = [double]"0.003530661738359763"

$ws.Range("A28").Value = [double]"239.9650551615946"
$ws.Range("B28").Value = [double]"248.628928923254"
$ws.Range("C28").Value = [double]"19.806366510482594"
$ws.Range("D28").Value = [double]"0.9690751318181683"
$ws.Range("E28").Value = [double]"0.2460274895046105"
$ws.Range("F28").Value = [double]"-0.005376053986508229"
$ws.Range("G28").Value = [double]"0.005194861952294624"

$ws.Range("A29").Value = [double]"239.9550278021346"
$ws.Range("B29").Value = [double]"248.51752949076078"
$ws.Range("C29").Value = [double]"19.927601523060158"
$ws.Range("D29").Value = [double]"0.9687457102920102"
$ws.Range("E29").Value = [double]"0.2441130627162939"
$ws.Range("F29").Value = [double]"-0.02199066840730142"
$ws.Range("G29").Value = [double]"0.03361564871786871"

$ws.Range("A30").Value = [double]"239.95457299609689"
$ws.Range("B30").Value = [double]"248.41190781593113"
$ws.Range("C30").Value = [double]"19.995327679843268"
$ws.Range("D30").Value = [double]"0.9703122025327043"
$ws.Range("E30").Value = [double]"0.24125401975862013"
$ws.Range("F30").Value = [double]"-0.0010395615503997334"
$ws.Range("G30").Value = [double]"0.0012889795837436196"

$ws.Range("A31").Value = [double]"239.95440941707568"
$ws.Range("B31").Value = [double]"248.31365248826282"
$ws.Range("C31").Value = [double]"19.999999954803343"
$ws.Range("D31").Value = [double]"0.9715983081968341"
$ws.Range("E31").Value = [double]"0.23607766896407623"
$ws.Range("F31").Value = [double]"-0.00039323171633810625"
$ws.Range("G31").Value = [double]"0.0003349302815016833"

$ws.Range("A32").Value = [double]"239.95426024858577"
$ws.Range("B32").Value = [double]"248.21576733353587"
$ws.Range("C32").Value = [double]"19.999957353639928"
$ws.Range("D32").Value = [double]"0.9715947774807931"
$ws.Range("E32").Value = [double]"0.23609215626350832"
$ws.Range("F32").Value = [double]"-0.000359966838787851"
$ws.Range("G32").Value = [double]"0.0003046354137425542"

$ws.Range("A33").Value = [double]"239.95412363572842"
$ws.Range("B33").Value = [double]"248.1178789797612"
$ws.Range("C33").Value = [double]"19.99995735998034"
$ws.Range("D33").Value = [double]"0.9715949117455783"
$ws.Range("E33").Value = [double]"0.23609168509871795"
$ws.Range("F33").Value = [double]"-0.00032965672705788054"
$ws.Range("G33").Value = [double]"0.0002790161981062121"

$ws.Range("A34").Value = [double]"239.95399949741363"
$ws.Range("B34").Value = [double]"248.0199907142764"
$ws.Range("C34").Value = [double]"19.99995736043383"
$ws.Range("D34").Value = [double]"0.9715949785697746"
$ws.Range("E34").Value = [double]"0.23609148173963343"
$ws.Range("F34").Value = [double]"-0.0002995547638872027"
$ws.Range("G34").Value = [double]"0.00025353770315696513"
